$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44358
$ws.Range("J2").Value = 200
$ws.Range("M2").Value = 650
$ws.Range("P2").Value = 108

# Row 3
$ws.Range("D3").Value = 44358
$ws.Range("J3").Value = 100

# Row 4
$ws.Range("D4").Value = 44188
$ws.Range("J4").Value = 200

# Row 5
$ws.Range("D5").Value = 44188
$ws.Range("J5").Value = 100

# Row 6
$ws.Range("D6").Value = 44321

# Row 7
$ws.Range("D7").Value = 44321

# Row 8
$ws.Range("D8").Value = 44308

# Row 9
$ws.Range("D9").Value = 44308

# Row 10
$ws.Range("D10").Value = 44335
$ws.Range("J10").Value = 150
$ws.Range("M10").Value = 633
$ws.Range("P10").Value = 106

# Row 11
$ws.Range("D11").Value = 44335
$ws.Range("J11").Value = 50

# Row 12
$ws.Range("D12").Value = 44230
$ws.Range("J12").Value = 100

# Row 13
$ws.Range("D13").Value = 44230
$ws.Range("J13").Value = 50

# Row 16
$ws.Range("D16").Value = 44554
$ws.Range("O16").Value = "Región de Ñuble"

# Row 17
$ws.Range("D17").Value = 44554
$ws.Range("O17").Value = "Región de Ñuble"

# Row 18
$ws.Range("D18").Value = 44491
$ws.Range("J18").Value = 200
$ws.Range("O18").Value = "Región Metropolitana"

# Row 19
$ws.Range("D19").Value = 44491
$ws.Range("J19").Value = 100
$ws.Range("O19").Value = "Región Metropolitana"

# Row 20
$ws.Range("D20").Value = 44293
$ws.Range("J20").Value = 100

# Row 21
$ws.Range("D21").Value = 44293
$ws.Range("J21").Value = 50

# Row 22
$ws.Range("D22").Value = 44525

# Row 23
$ws.Range("D23").Value = 44525
